# Regenerate the "K" column (column G) values on Sheet1.
# The source data now derives strikeouts (K) differently ("use K instead of
# Strike#"), so the per-appearance K values below replace the old ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    3  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 3
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    14 = 1
    15 = 0
    16 = 1
    17 = 0
    19 = 1
    20 = 1
    21 = 1
    22 = 2
    23 = 0
    24 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
